$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "short-url" column (B): rename shared value "pacLD1" -> "KYq76X" for all data rows
$ws.Range("B2:B4").Value = "KYq76X"

# "oip" column (U): rename shared value "null" -> "-" for all data rows
$ws.Range("U2:U4").Value = "-"

# "hst" column (V): value becomes "-" (like the oip column) and alignment switches
# from right-aligned to left-aligned to match the new text style
$ws.Range("V2:V4").Value = "-"
$ws.Range("V2:V4").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
